$wb = $excel.ActiveWorkbook

# Remove Sheet2 and Sheet3, keeping only Sheet1
$wb.Worksheets.Item("Sheet2").Delete() | Out-Null
$wb.Worksheets.Item("Sheet3").Delete() | Out-Null

$ws = $wb.Worksheets.Item("Sheet1")

# Update header row
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Value"

# Add new data rows
$ws.Range("A2").Value = "Item1"
$ws.Range("B2").Value = 100

$ws.Range("A3").Value = "Item2"
$ws.Range("B3").Value = 200
